{"js": "// Rewrite the intro paragraph of the \"Git\" document:\n//   \"Git est une s\u00e9rie de logiciels d\u00e9velopp\u00e9es par Linus Torvald aussi\n//    connu pour avoir d\u00e9velopp\u00e9 Linux, Git est une \"\n// becomes\n//   \"Git est un VCS cr\u00e9e par Linus Torvald (aussi connu pour avoir\n//    d\u00e9velopp\u00e9 Linux), son but est de vouloir faciliter les projets\n//    d'\u00e9quipes en donnant acc\u00e8s \u00e0 une arborescence mutualis\u00e9e pour\n//    permettre aux utilisateurs concern\u00e9s d'ajouter, de supprimer et de\n//    modifier des \u00e9l\u00e9ments dans ces dossiers partag\u00e9s. Le langage git\n//    (tout comme le langage shell) est un langage de script, cependant,\n//    ce langage n'a pas la possibilit\u00e9 d'ex\u00e9cuter ou de lire un fichier.\"\n//\n// The word \"Torvald\" sits in its own run (flanked by proofErr spell-check\n// markers) and must be left untouched, so the paragraph is edited in two\n// surgical pieces: the text before \"Torvald\" and the text after it. Using\n// Range.insertText(..., Word.InsertLocation.replace) on ranges that abut\n// existing runs lets Word inherit the surrounding (Times New Roman)\n// character formatting instead of falling back to run defaults.\n\nconst body = context.document.body;\n\n// --- Piece 1: \"Git est une s\u00e9rie ... par Linus \" -> \"Git est un VCS cr\u00e9e par Linus \"\nconst beforeTorvald = body.search(\n  \"Git est une s\u00e9rie de logiciels d\u00e9velopp\u00e9es par Linus \",\n  { matchCase: true }\n);\nbeforeTorvald.load(\"text\");\nawait context.sync();\n\nif (beforeTorvald.items.length > 0) {\n  beforeTorvald.items[0].insertText(\n    \"Git est un VCS cr\u00e9e par Linus \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Piece 2: \" aussi connu ... Git est une \" -> \" (aussi connu ... fichier.\"\nconst afterTorvald = body.search(\n  \" aussi connu pour avoir d\u00e9velopp\u00e9 Linux, Git est une \",\n  { matchCase: true }\n);\nafterTorvald.load(\"text\");\nawait context.sync();\n\nif (afterTorvald.items.length > 0) {\n  const replacement =\n    \" (aussi connu pour avoir d\u00e9velopp\u00e9 Linux), son but est de vouloir \" +\n    \"faciliter les projets d\\u2019\u00e9quipes en donnant acc\u00e8s \u00e0 une \" +\n    \"arborescence mutualis\u00e9e pour permettre aux utilisateurs concern\u00e9s \" +\n    \"d\\u2019ajouter, de supprimer et de modifier des \u00e9l\u00e9ments dans ces \" +\n    \"dossiers partag\u00e9s. Le langage git (tout comme le langage shell) est \" +\n    \"un langage de script, cependant, ce langage n\\u2019a pas la \" +\n    \"possibilit\u00e9 d\\u2019ex\u00e9cuter ou de lire un fichier.\";\n\n  afterTorvald.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Rewrite the intro paragraph of the \"Git\" document:\n#   \"Git est une s\u00e9rie de logiciels d\u00e9velopp\u00e9es par Linus Torvald aussi\n#    connu pour avoir d\u00e9velopp\u00e9 Linux, Git est une \"\n# becomes\n#   \"Git est un VCS cr\u00e9e par Linus Torvald (aussi connu pour avoir\n#    d\u00e9velopp\u00e9 Linux), son but est de vouloir faciliter les projets\n#    d'\u00e9quipes en donnant acc\u00e8s \u00e0 une arborescence mutualis\u00e9e pour\n#    permettre aux utilisateurs concern\u00e9s d'ajouter, de supprimer et de\n#    modifier des \u00e9l\u00e9ments dans ces dossiers partag\u00e9s. Le langage git\n#    (tout comme le langage shell) est un langage de script, cependant,\n#    ce langage n'a pas la possibilit\u00e9 d'ex\u00e9cuter ou de lire un fichier.\"\n#\n# \"Torvald\" sits in its own run flanked by proofErr spell-check markers\n# and must be left untouched, so the paragraph is edited in two pieces:\n# the text before \"Torvald\" and the text after it. Assigning .Text on a\n# Find-located Range (rather than deleting + re-inserting) lets Word\n# keep the surrounding (Times New Roman) character formatting.\n\n$d = $word.ActiveDocument\n\n# --- Piece 1: \"Git est une s\u00e9rie ... par Linus \" -> \"Git est un VCS cr\u00e9e par Linus \"\n$rng1 = $d.Content\nif ($rng1.Find.Execute(\"Git est une s\u00e9rie de logiciels d\u00e9velopp\u00e9es par Linus \")) {\n    $rng1.Text = \"Git est un VCS cr\u00e9e par Linus \"\n}\n\n# --- Piece 2: \" aussi connu ... Git est une \" -> \" (aussi connu ... fichier.\"\n$replacement = \" (aussi connu pour avoir d\u00e9velopp\u00e9 Linux), son but est de vouloir faciliter les projets d\u2019\u00e9quipes en donnant acc\u00e8s \u00e0 une arborescence mutualis\u00e9e pour permettre aux utilisateurs concern\u00e9s d\u2019ajouter, de supprimer et de modifier des \u00e9l\u00e9ments dans ces dossiers partag\u00e9s. Le langage git (tout comme le langage shell) est un langage de script, cependant, ce langage n\u2019a pas la possibilit\u00e9 d\u2019ex\u00e9cuter ou de lire un fichier.\"\n\n$rng2 = $d.Content\nif ($rng2.Find.Execute(\" aussi connu pour avoir d\u00e9velopp\u00e9 Linux, Git est une \")) {\n    $rng2.Text = $replacement\n}\n"}
